$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("partner")
Write-Host $ws.Name
